# Add two more datetime rows (4 and 5), mirroring the pattern already
# used in rows 2-3: col A holds a raw Excel date serial (formatted via
# the existing date/time style), col B holds a TEXT() formula that
# renders it as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date serials for 2017-01-02 11:30:00 and 2017-01-03 11:30:00
$ws.Range("A4").Value = 42737.479166666664
$ws.Range("A5").Value = 42738.479166666664

# Reuse the same number format/style as the existing date cells (A2)
# by copying its formatting onto the new cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill B4:B5 with the same TEXT() formula pattern as B2/B3; Excel will
# store this as a shared formula across the B4:B5 block.
$ws.Range("B4:B5").Formula = '=TEXT(A4,"yyyy-mm-dd hh:mm:ss")'

# Match the final selection state recorded in the sheet.
$ws.Range("A2").Select() | Out-Null
